$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.648.72"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "3.500.35"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "'169.29"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +4.93%  "
$ws.Range("D9").Value = "'0.130"
$ws.Range("E9").Value = "  +5.92%  "
$ws.Range("D10").Value = "'7.33"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").Value = "'0.432"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "4.103.37"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "'28.23"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").Value = "'0.0000180"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "66.654.53"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").Value = "3.508.75"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "'6.33"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").Value = "'394.44"
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").Value = "'7.95"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "'73.21"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").Value = "'0.534"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("D25").Value = "'0.0000122"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").Value = "'10.20"
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'6.34"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "'1.46"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").Value = "'2.06"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").Value = "'23.81"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("D33").Value = "'7.37"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("E34").Value = "  +4.95%  "
$ws.Range("D35").Value = "'162.12"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").Value = "'0.901"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("E37").Value = "  +2.49%  "
$ws.Range("E38").Value = "  +2.12%  "
$ws.Range("D39").Value = "'4.66"
$ws.Range("E39").Value = "  +4.13%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0742"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "'26.48"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'26.98"
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.801.71"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "'42.91"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").Value = "'2.57"
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("D46").Value = "'0.0312"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").Value = "'342.86"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").Value = "'33.99"
$ws.Range("E49").Value = "  +4.77%  "
$ws.Range("D50").Value = "'0.855"
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("D51").Value = "'6.52"
$ws.Range("E51").Value = "  +1.94%  "
